$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text update (row 22, column B: resistor designators) ---
$ws.Range("B22").Value = "R17,R18,R1,R8"

# --- Fill in / renumber the BOM line-item numbers in column A ---
# Row 22 previously had no line number; it now becomes item 21.
$ws.Range("A22").Value = 21
# Row 23 shifts from 21 to 22.
$ws.Range("A23").Value = 22
# Rows 24 and 25 previously had no line number; they become 23 and 24.
$ws.Range("A24").Value = 23
$ws.Range("A25").Value = 24
# Rows 26-32 shift up by three (22->25 ... 28->31).
$ws.Range("A26").Value = 25
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27
$ws.Range("A29").Value = 28
$ws.Range("A30").Value = 29
$ws.Range("A31").Value = 30
$ws.Range("A32").Value = 31
# Row 33 no longer carries a line number.
$ws.Range("A33").ClearContents() | Out-Null

# --- Quantity changes (column D) that ripple into the Extended Price formula (column I) ---
# Row 22: quantity 2 -> 4 (I22 = H22*D22 recalculates automatically).
$ws.Range("D22").Value = 4
# Row 31: quantity 4 -> 2 (I31 = H31*D31 recalculates automatically).
$ws.Range("D31").Value = 2

# Row 33 reverts to the sheet's default (non-custom) row height.
$ws.Rows.Item(33).AutoFit() | Out-Null

# --- View state: scroll back to the left edge and move the active selection ---
$ws.Activate() | Out-Null
$ws.Range("E37").Select() | Out-Null
